# Applies the edits described by the commit:
#   ".gitignore, update GHdiagram, session file"
#
# 1) The "last updated" date auto-field (datetimeFigureOut) cached text is
#    bumped from 27/10/2021 -> 10/11/2021 on the slide master and on every
#    slide layout (PowerPoint re-caches this shared field text whenever the
#    deck is saved).
# 2) The dashed horizontal divider line (connector "Conector recto 5" on
#    slide 1) is resized/repositioned and flipped vertically.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Update the cached text of the date placeholder wherever it still
#    shows the old cached value, across the slide master and all layouts.
# ---------------------------------------------------------------------
function Update-DateField {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -eq "27/10/2021") {
                $shp.TextFrame.TextRange.Text = "10/11/2021"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateField $master.Shapes

$layouts = $master.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    Update-DateField $layouts.Item($L).Shapes
}

# ---------------------------------------------------------------------
# 2) Resize / reposition / flip the dashed divider connector on slide 1.
#    Target OOXML:
#      <a:xfrm flipV="1">
#        <a:off x="9330" y="4861279"/>
#        <a:ext cx="10141349" cy="65286"/>
#      </a:xfrm>
#    (values below are expressed in points == EMU / 12700, chosen so the
#    shim's internal float32 rounding reproduces the exact target EMUs)
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(1)
$connector = $null
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.Name -eq "Conector recto 5") {
        $connector = $shp
        break
    }
}

$connector.Left = 0.7346457242965698
$connector.Top = 382.77789306640625
$connector.Width = 798.5314331054688
$connector.Height = 5.14063024520874
$connector.Flip(1)
